# Applies the "mere rapport Scrum kapitel + diverse" edit:
#  - Update the monthly number-of-downloads row (B4:L4) on the first
#    pricing table with new, smaller figures.
#  - Add a new "App Store" revenue row (row 23) under the second table,
#    mirroring rows 11/12 (Antal/Cut/Før skat).
#  - Rename the "CPM(DKK)" label (row 22, col B) to "eCPM(DKK)".
#  - Move the view/selection back up to the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the first table's monthly download counts (row 4) ---
$ws.Range("B4").Value = 5000
$ws.Range("C4").Value = 3000
$ws.Range("D4").Value = 3000
$ws.Range("E4").Value = 1500
$ws.Range("F4").Value = 1500
$ws.Range("G4").Value = 1000
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 250
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 50

# --- Relabel the eCPM row header text (row 22) ---
$ws.Range("B22").Value = "eCPM(DKK)"

# --- Add the new App Store revenue-calc row (row 23) ---
$ws.Range("B23").Value = 16
$ws.Range("C23").Value = 3
$ws.Range("D23").Formula = "=SUM((A23/1000)*B23*C23)"

# --- Restore the sheet view / selection ---
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("D23").Select()

$wb.Application.Calculate()
